$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1442.8889
$ws.Range("I40").Value = 1400.1538
$ws.Range("K40").Value = 1400.1538
$ws.Range("M40").Value = -1225.1538
# Row 51
$ws.Range("H51").Value = 2640.6155
$ws.Range("J51").Value = 2884.6667
$ws.Range("L51").Value = 2884.6667
$ws.Range("N51").Value = -3852.6667
# Row 58
$ws.Range("H58").Value = 1143.8889
$ws.Range("I58").Value = 233.33333
$ws.Range("J58").Value = 1326
$ws.Range("K58").Value = 699.99999
$ws.Range("L58").Value = 3978
$ws.Range("M58").Value = -549.99999
$ws.Range("N58").Value = -4278
# Row 98
$ws.Range("H98").Value = 3442362.5
$ws.Range("I98").Value = 4516992
$ws.Range("J98").Value = 3548.1
$ws.Range("K98").Value = 4516992
$ws.Range("L98").Value = 3548.1
$ws.Range("M98").Value = -4515494
$ws.Range("N98").Value = -6544.1
# Row 107
$ws.Range("H107").Value = 640.3333
$ws.Range("I107").Value = 292.33334
$ws.Range("J107").Value = 1858.3334
$ws.Range("K107").Value = 292.33334
$ws.Range("L107").Value = 1858.3334
$ws.Range("M107").Value = 1627.66666
$ws.Range("N107").Value = -5698.3334
# Row 112
$ws.Range("H112").Value = 6499.263
$ws.Range("J112").Value = 6499.263
$ws.Range("L112").Value = 19497.789
$ws.Range("N112").Value = -21713.789
# Row 122
$ws.Range("H122").Value = 3442362.5
$ws.Range("I122").Value = 4516992
$ws.Range("J122").Value = 3548.1
$ws.Range("K122").Value = 13550976
$ws.Range("L122").Value = 10644.3
$ws.Range("M122").Value = -13548526
$ws.Range("N122").Value = -15544.3
# Row 132
$ws.Range("H132").Value = 4349079
$ws.Range("I132").Value = 4349079
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13047237
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -13044707
$ws.Range("N132").Value = $null
# Row 137
$ws.Range("H137").Value = 960.4167
$ws.Range("I137").Value = 927
$ws.Range("J137").Value = 1036.3636
$ws.Range("K137").Value = 2781
$ws.Range("L137").Value = 3109.0908
$ws.Range("M137").Value = -231
$ws.Range("N137").Value = -8209.0908

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 296403.53
$ws.Range("I32").Value = 2860.557
$ws.Range("J32").Value = 1516924.4
$ws.Range("K32").Value = 2860.557
$ws.Range("L32").Value = 1516924.4
$ws.Range("M32").Value = -2573.557
$ws.Range("N32").Value = -1517498.4
# Row 52
$ws.Range("H52").Value = 29746.666
$ws.Range("J52").Value = 29746.666
$ws.Range("L52").Value = 29746.666
$ws.Range("N52").Value = -30382.666
# Row 61
$ws.Range("H61").Value = 1869.7906
$ws.Range("I61").Value = 1546.1515
$ws.Range("K61").Value = 1546.1515
$ws.Range("M61").Value = -1334.1515
# Row 136
$ws.Range("H136").Value = 1869.7906
$ws.Range("I136").Value = 1546.1515
$ws.Range("K136").Value = 4638.4545
$ws.Range("M136").Value = -2088.4545

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3061.4285
$ws.Range("I105").Value = 2681.6667
$ws.Range("J105").Value = 3567.7778
$ws.Range("K105").Value = 2681.6667
$ws.Range("L105").Value = 3567.7778
$ws.Range("M105").Value = -934.6667000000002
$ws.Range("N105").Value = -7061.7778
# Row 134
$ws.Range("H134").Value = 2350.25
$ws.Range("I134").Value = 1921.3611
$ws.Range("J134").Value = 4280.25
$ws.Range("K134").Value = 5764.0833
$ws.Range("L134").Value = 12840.75
$ws.Range("M134").Value = -3229.0833
$ws.Range("N134").Value = -17910.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 75.916664
$ws.Range("I7").Value = 61.285713
$ws.Range("J7").Value = 96.40000000000001
$ws.Range("K7").Value = 61.285713
$ws.Range("L7").Value = 96.40000000000001
$ws.Range("M7").Value = 51.714287
$ws.Range("N7").Value = -322.4
# Row 31
$ws.Range("H31").Value = 1549625.4
$ws.Range("I31").Value = 2943372.2
$ws.Range("J31").Value = 1017.6667
$ws.Range("K31").Value = 2943372.2
$ws.Range("L31").Value = 1017.6667
$ws.Range("M31").Value = -2943077.2
$ws.Range("N31").Value = -1607.6667
# Row 34
$ws.Range("H34").Value = 1549625.4
$ws.Range("I34").Value = 2943372.2
$ws.Range("J34").Value = 1017.6667
$ws.Range("K34").Value = 2943372.2
$ws.Range("L34").Value = 1017.6667
$ws.Range("M34").Value = -2943170.2
$ws.Range("N34").Value = -1421.6667
# Row 58
$ws.Range("H58").Value = 795.11365
$ws.Range("I58").Value = 687.44116
$ws.Range("J58").Value = 1161.2
$ws.Range("K58").Value = 687.44116
$ws.Range("L58").Value = 1161.2
$ws.Range("M58").Value = -484.44116
$ws.Range("N58").Value = -1567.2
# Row 122
$ws.Range("H122").Value = 5714867
$ws.Range("I122").Value = 586.2174
$ws.Range("J122").Value = 16667238
$ws.Range("K122").Value = 1758.6522
$ws.Range("L122").Value = 50001714
$ws.Range("M122").Value = 691.3478
$ws.Range("N122").Value = -50006614
# Row 132
$ws.Range("H132").Value = 4168928.2
$ws.Range("I132").Value = 1751.1052
$ws.Range("J132").Value = 20004202
$ws.Range("K132").Value = 5253.3156
$ws.Range("L132").Value = 60012606
$ws.Range("M132").Value = -2723.3156
$ws.Range("N132").Value = -60017666
# Row 134
$ws.Range("H134").Value = 3682.1365
$ws.Range("I134").Value = 4122.6284
$ws.Range("J134").Value = 1969.1111
$ws.Range("K134").Value = 12367.8852
$ws.Range("L134").Value = 5907.3333
$ws.Range("M134").Value = -9832.885199999999
$ws.Range("N134").Value = -10977.3333
# Row 136
$ws.Range("H136").Value = 795.11365
$ws.Range("I136").Value = 687.44116
$ws.Range("J136").Value = 1161.2
$ws.Range("K136").Value = 2062.32348
$ws.Range("L136").Value = 3483.6
$ws.Range("M136").Value = 487.67652
$ws.Range("N136").Value = -8583.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 951.2973
$ws.Range("I68").Value = 680.55554
$ws.Range("J68").Value = 1207.7894
$ws.Range("K68").Value = 2041.66662
$ws.Range("L68").Value = 3623.3682
$ws.Range("M68").Value = -1230.66662
$ws.Range("N68").Value = -5245.3682
# Row 71
$ws.Range("H71").Value = 951.2973
$ws.Range("I71").Value = 680.55554
$ws.Range("J71").Value = 1207.7894
$ws.Range("K71").Value = 6124.99986
$ws.Range("L71").Value = 10870.1046
$ws.Range("M71").Value = -2068.99986
$ws.Range("N71").Value = -18982.1046
# Row 74
$ws.Range("H74").Value = 6887.5713
$ws.Range("J74").Value = 8000
$ws.Range("L74").Value = 24000
$ws.Range("N74").Value = -26122
# Row 77
$ws.Range("H77").Value = 6887.5713
$ws.Range("J77").Value = 8000
$ws.Range("L77").Value = 72000
$ws.Range("N77").Value = -82608
# Row 102
$ws.Range("H102").Value = 5846.3
$ws.Range("I102").Value = 2926
$ws.Range("K102").Value = 8778
$ws.Range("M102").Value = -6344
# Row 105
$ws.Range("H105").Value = 1764234
$ws.Range("J105").Value = 1860858.1
$ws.Range("L105").Value = 5582574.300000001
$ws.Range("N105").Value = -5587816.300000001
# Row 113
$ws.Range("H113").Value = 892.7324
$ws.Range("I113").Value = 686.1
$ws.Range("J113").Value = 926.60657
$ws.Range("K113").Value = 2058.3
$ws.Range("L113").Value = 2779.81971
$ws.Range("M113").Value = 111.6999999999998
$ws.Range("N113").Value = -7119.81971
# Row 134
$ws.Range("H134").Value = 2826.25
$ws.Range("I134").Value = 2230
$ws.Range("K134").Value = 6690
$ws.Range("M134").Value = -1620

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 19588.582
$ws.Range("I132").Value = 1092.0682
$ws.Range("J132").Value = 93574.63
$ws.Range("K132").Value = 3276.2046
$ws.Range("L132").Value = 280723.89
$ws.Range("M132").Value = -746.2046
$ws.Range("N132").Value = -285783.89

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 25
$ws.Range("H25").Value = 16666.666
$ws.Range("I25").Value = 16666.666
$ws.Range("K25").Value = 16666.666
$ws.Range("M25").Value = -16436.666
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null
# Row 100
$ws.Range("H100").Value = 2176.7585
$ws.Range("I100").Value = 1183.4166
$ws.Range("J100").Value = 2877.9412
$ws.Range("K100").Value = 1183.4166
$ws.Range("L100").Value = 2877.9412
$ws.Range("M100").Value = -642.4166
$ws.Range("N100").Value = -3959.9412
# Row 122
$ws.Range("H122").Value = 50003176
$ws.Range("I122").Value = 2577.8333
$ws.Range("J122").Value = 125004070
$ws.Range("K122").Value = 7733.499899999999
$ws.Range("L122").Value = 375012210
$ws.Range("M122").Value = -5283.499899999999
$ws.Range("N122").Value = -375017110
# Row 132
$ws.Range("H132").Value = 2382642.8
$ws.Range("I132").Value = 3247883.8
$ws.Range("J132").Value = 3230.25
$ws.Range("K132").Value = 9743651.399999999
$ws.Range("L132").Value = 9690.75
$ws.Range("M132").Value = -9741121.399999999
$ws.Range("N132").Value = -14750.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 4999.3335
$ws.Range("J40").Value = 4999.3335
$ws.Range("L40").Value = 4999.3335
$ws.Range("N40").Value = -5297.3335
# Row 126
$ws.Range("H126").Value = 2012.8572
$ws.Range("J126").Value = 4000
$ws.Range("L126").Value = 12000
$ws.Range("N126").Value = -16940
# Row 132
$ws.Range("H132").Value = 26599912
$ws.Range("I132").Value = 36765944
$ws.Range("J132").Value = 11826.385
$ws.Range("K132").Value = 110297832
$ws.Range("L132").Value = 35479.155
$ws.Range("M132").Value = -110295302
$ws.Range("N132").Value = -40539.155
# Row 136
$ws.Range("H136").Value = 67012.13
$ws.Range("I136").Value = 71727.28999999999
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 215181.87
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -212631.87
$ws.Range("N136").Value = -8100
